$d = $word.ActiveDocument

# The sentence "Корисник притиска дугме „Даље“." becomes
# "Корисник притиска дугме за наставак." - i.e. the quoted
# button-name „Даље“ is replaced by "за наставак", ending up
# split across two runs while the trailing ". " run is left alone.

# 1) Replace the quoted button name with the new wording.
$rng = $d.Content
$rng.Find.Execute("„Даље“", $true, $false, $false, $false, $false, $true, 1, $false, "за наставак", 2)

# 2) Force a run boundary right before "за наставак" so it ends up
#    as its own run (mirrors the original author's edit, which left
#    "Корисник притиска дугме " as one run and "за наставак" as another).
$wedge1 = $d.Content
$wedge1.Find.Execute("притиска дугме ")
$p1 = $wedge1.End
$bm1 = $d.Range($p1, $p1)
$d.Bookmarks.Add("TmpSplit1", $bm1)
$d.Bookmarks("TmpSplit1").Delete()

# 3) Force a run boundary right after "за наставак" so the following
#    ". " text remains its own, untouched run.
$wedge2 = $d.Content
$wedge2.Find.Execute("притиска дугме за наставак")
$p2 = $wedge2.End
$bm2 = $d.Range($p2, $p2)
$d.Bookmarks.Add("TmpSplit2", $bm2)
$d.Bookmarks("TmpSplit2").Delete()
